# Add new Testcase in MyLocationTest Class
# Appends a new key/value row (locationseach / fgh) to the test data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = "locationseach"
$ws.Range("B27").Value = "fgh"

# Bring the new row into view / select it, like Excel would leave it after
# the user typed the new data in at the bottom of the sheet.
$ws.Range("B27").Select()
